$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 42, shifting existing rows 42-143 down to 43-144.
$ws.Rows.Item(42).Insert()

# Populate the newly inserted row 42 with the new data point.
$ws.Cells.Item(42, 4).Value = 44979
$ws.Cells.Item(42, 11).Value = "Thompson seedless"
$ws.Cells.Item(42, 13).Value = 400
$ws.Cells.Item(42, 14).Value = 10000
$ws.Cells.Item(42, 15).Value = 11000
$ws.Cells.Item(42, 16).Value = 10500
$ws.Cells.Item(42, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(42, 19).Value = 583

# Columns that remain identical to the surrounding rows (same as what Insert
# carried down from row 41 is NOT guaranteed, so set them explicitly).
$ws.Cells.Item(42, 1).Value = 8
$ws.Cells.Item(42, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(42, 3).Value = "Coquimbo"
$ws.Cells.Item(42, 5).Value = 4
$ws.Cells.Item(42, 6).Value = "Fruta"
$ws.Cells.Item(42, 7).Value = 100109
$ws.Cells.Item(42, 8).Value = "Uva"
$ws.Cells.Item(42, 9).Value = 100109001
$ws.Cells.Item(42, 10).Value = "Uva"
$ws.Cells.Item(42, 12).Value = "Primera"
$ws.Cells.Item(42, 17).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(42, 20).Value = 18
